$d = $word.ActiveDocument
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Find-ParagraphIndex($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# --- Part 1 ---------------------------------------------------------
# Collapse the entire model-summary table block (the box-drawing table,
# layer rows, and the "Total/Trainable/Non-trainable/Optimizer params"
# lines, plus the trailing page-break paragraph) together with the
# "Hi, after thinking..." paragraph that follows it into two plain
# paragraphs: a short "From Santi:" note (still in Consolas, carrying
# the page-break render hint) followed by the original "Hi, after
# thinking..." text with all special formatting stripped.

$tableStartIdx = Find-ParagraphIndex $d "┌─────"
$helloIdx = Find-ParagraphIndex $d "Hi, after thinking more about it"

if ($tableStartIdx -eq -1 -or $helloIdx -eq -1) {
    throw "Could not locate model-summary table / intro paragraph"
}

$pStart = $d.Paragraphs.Item($tableStartIdx)
$pEnd = $d.Paragraphs.Item($helloIdx)
$range1 = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$xml1 = "<w:p $ns>" +
            "<w:pPr><w:rPr><w:rFonts w:ascii='Consolas' w:hAnsi='Consolas'/></w:rPr></w:pPr>" +
            "<w:r>" +
                "<w:rPr><w:rFonts w:ascii='Consolas' w:hAnsi='Consolas'/></w:rPr>" +
                "<w:lastRenderedPageBreak/>" +
                "<w:t>From Santi:</w:t>" +
            "</w:r>" +
        "</w:p>" +
        "<w:p $ns>" +
            "<w:r>" +
                "<w:t>Hi, after thinking more about it, I realized that it would be more interesting to analyze the performance of the model with the night view. A common problem we faced when analyzing the photos is that</w:t>" +
            "</w:r>" +
        "</w:p>"

$range1.InsertXML($xml1)

# --- Part 2 ---------------------------------------------------------
# The page-break render hint that used to sit on the "New labeled
# dataset" bullet now belongs on the "When I mention Storm ###"
# paragraph just above it.

$d2 = $word.ActiveDocument
$stormIdx = Find-ParagraphIndex $d2 "When I mention Storm ###"
$labeledIdx = Find-ParagraphIndex $d2 "New labeled dataset"

if ($stormIdx -eq -1 -or $labeledIdx -eq -1) {
    throw "Could not locate Storm ### / New labeled dataset paragraphs"
}

$pStorm = $d2.Paragraphs.Item($stormIdx)
$xmlStorm = "<w:p $ns><w:r><w:lastRenderedPageBreak/><w:t xml:space='preserve'>When I mention Storm ### I mean the number of Column C (`"Storm Number`") of the Excel spreadsheet </w:t></w:r></w:p>"
$pStorm.Range.InsertXML($xmlStorm)

$d3 = $word.ActiveDocument
$labeledIdx = Find-ParagraphIndex $d3 "New labeled dataset"
$pLabeled = $d3.Paragraphs.Item($labeledIdx)
$xmlLabeled = "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr></w:pPr><w:r><w:t>New labeled dataset</w:t></w:r></w:p>"
$pLabeled.Range.InsertXML($xmlLabeled)

Write-Output "done"
